# Add two new "problems solved" entries to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 45: Cute little cat (uses the existing "Trie" classname entry)
$ws.Cells.Item(45, 1).Value = 48
$ws.Cells.Item(45, 2).Value = "Cute little cat"
$ws.Cells.Item(45, 3).Value = "Trie"
$ws.Cells.Item(45, 8).Value = "CuteCatImpl"

# Row 46: BiggestXOR problem
$ws.Cells.Item(46, 1).Value = 49
$ws.Cells.Item(46, 2).Value = "BiggestXOR problem"
$ws.Cells.Item(46, 8).Value = "BiggestXOR"

# Match the selection left behind by the edit.
$ws.Range("H46").Select()
